# Insert a new weekly record at the top of the Brócoli price series for
# "Feria Lagunitas de Puerto Montt". This shifts every existing data row
# (174:239) down by one row (175:240) and fills the freshly opened row 174
# with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 174:239 down to 175:240, opening up a blank row 174.
$ws.Rows("174:174").Insert()

# Populate the new row 174 with the latest price observation.
$ws.Range("A174").Value = 4
$ws.Range("B174").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C174").Value = "Los Lagos"
$ws.Range("D174").Value = 44524
$ws.Range("E174").Value = 10
$ws.Range("F174").Value = 100112023
$ws.Range("G174").Value = "Brócoli"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Segunda"
$ws.Range("J174").Value = 100
$ws.Range("K174").Value = 1000
$ws.Range("L174").Value = 1000
$ws.Range("M174").Value = 1000
$ws.Range("N174").Value = "$/unidad"
$ws.Range("O174").Value = "Región Metropolitana"
$ws.Range("P174").Value = 1000
$ws.Range("Q174").Value = 1
$ws.Range("R174").Value = "Hortaliza"
